$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 201 (the old last data row), shifting it down to row 203
$ws.Rows.Item(201).Resize(2).Insert()

# New row 201: Choclo / Choclero / Primera, date 2022-04-05 (44656)
$ws.Range("A201").Value = 7
$ws.Range("B201").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C201").Value = "Ñuble"
$ws.Range("D201").Value = 44656
$ws.Range("D201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E201").Value = 16
$ws.Range("F201").Value = 100112024
$ws.Range("G201").Value = "Choclo"
$ws.Range("H201").Value = "Choclero"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 10000
$ws.Range("K201").Value = 250
$ws.Range("L201").Value = 250
$ws.Range("M201").Value = 250
$ws.Range("N201").Value = "`$/unidad"
$ws.Range("O201").Value = "Región del Maule"
$ws.Range("P201").Value = 250
$ws.Range("Q201").Value = 1
$ws.Range("R201").Value = "Hortaliza"

# New row 202: Choclo / Choclero / Segunda, date 2022-04-05 (44656)
$ws.Range("A202").Value = 7
$ws.Range("B202").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C202").Value = "Ñuble"
$ws.Range("D202").Value = 44656
$ws.Range("D202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E202").Value = 16
$ws.Range("F202").Value = 100112024
$ws.Range("G202").Value = "Choclo"
$ws.Range("H202").Value = "Choclero"
$ws.Range("I202").Value = "Segunda"
$ws.Range("J202").Value = 8000
$ws.Range("K202").Value = 180
$ws.Range("L202").Value = 180
$ws.Range("M202").Value = 180
$ws.Range("N202").Value = "`$/unidad"
$ws.Range("O202").Value = "Región del Maule"
$ws.Range("P202").Value = 180
$ws.Range("Q202").Value = 1
$ws.Range("R202").Value = "Hortaliza"
